$wb = $excel.ActiveWorkbook

# ---- Sheet: P_valores ----
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.07634306692504245
$wsP.Range("D2").Value = 0.3920087038318165
$wsP.Range("E2").Value = 0.3918716475723727
$wsP.Range("F2").Value = 0.09458195244530843

$wsP.Range("B3").Value = 0.07634306692504245
$wsP.Range("D3").Value = 0.4295151170385405
$wsP.Range("E3").Value = 0.1656384087488514
$wsP.Range("F3").Value = 0.4619385745363145

$wsP.Range("B4").Value = 0.3920087038318165
$wsP.Range("C4").Value = 0.4295151170385405
$wsP.Range("E4").Value = 0.5972245527549576
$wsP.Range("F4").Value = 0.6749747880229071

$wsP.Range("B5").Value = 0.3918716475723727
$wsP.Range("C5").Value = 0.1656384087488514
$wsP.Range("D5").Value = 0.5972245527549576
$wsP.Range("F5").Value = 0.3226275213869265

$wsP.Range("B6").Value = 0.09458195244530843
$wsP.Range("C6").Value = 0.4619385745363145
$wsP.Range("D6").Value = 0.6749747880229071
$wsP.Range("E6").Value = 0.3226275213869265

# ---- Sheet: Estadisticos_DM ----
$wsD = $wb.Worksheets.Item("Estadisticos_DM")

$wsD.Range("C2").Value = -1.859765880946411
$wsD.Range("D2").Value = -0.8731568255455207
$wsD.Range("E2").Value = -0.8734139177567816
$wsD.Range("F2").Value = -1.747017215646653

$wsD.Range("B3").Value = 1.859765880946411
$wsD.Range("D3").Value = 0.8048569576179185
$wsD.Range("E3").Value = 1.433984059962139
$wsD.Range("F3").Value = 0.7487502464899295

$wsD.Range("B4").Value = 0.8731568255455207
$wsD.Range("C4").Value = -0.8048569576179185
$wsD.Range("E4").Value = 0.5361651101071694
$wsD.Range("F4").Value = -0.4249859301041415

$wsD.Range("B5").Value = 0.8734139177567816
$wsD.Range("C5").Value = -1.433984059962139
$wsD.Range("D5").Value = -0.5361651101071694
$wsD.Range("F5").Value = -1.011811803301595

$wsD.Range("B6").Value = 1.747017215646653
$wsD.Range("C6").Value = -0.7487502464899295
$wsD.Range("D6").Value = 0.4249859301041415
$wsD.Range("E6").Value = 1.011811803301595
